# Update localization status report: change Status from "Ready for handoff"
# to "In Translation" for the two files whose latest handoff is still in
# flight (10bb3719-...md and 66575674-...md), on all three sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns B (zh-cn) and C (de-de) hold the per-language status ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "In Translation"
$wsOverview.Range("C3").Value = "In Translation"
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"

# --- zh-cn detail sheet: column C holds Status ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C4").Value = "In Translation"

# --- de-de detail sheet: column C holds Status ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C4").Value = "In Translation"
